$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 56.555557
$ws.Range("I5").Value = 44.142857
$ws.Range("K5").Value = 44.142857
$ws.Range("M5").Value = 70.85714300000001

$ws.Range("H9").Value = 427.58334
$ws.Range("I9").Value = 367.375
$ws.Range("J9").Value = 548
$ws.Range("K9").Value = 367.375
$ws.Range("L9").Value = 548
$ws.Range("M9").Value = -198.375
$ws.Range("N9").Value = -886

$ws.Range("H11").Value = 50.7
$ws.Range("I11").Value = 50.7
$ws.Range("K11").Value = 50.7
$ws.Range("M11").Value = 89.3

$ws.Range("H12").Value = 1381.8182
$ws.Range("I12").Value = 171.57143
$ws.Range("K12").Value = 171.57143
$ws.Range("M12").Value = -1.571429999999992

$ws.Range("H53").Value = 1398.6471
$ws.Range("I53").Value = 150.8
$ws.Range("J53").Value = 1918.5834
$ws.Range("K53").Value = 150.8
$ws.Range("L53").Value = 1918.5834
$ws.Range("M53").Value = 486.2
$ws.Range("N53").Value = -3192.5834

$ws.Range("H54").Value = 501000
$ws.Range("I54").Value = 501000
$ws.Range("K54").Value = 501000
$ws.Range("M54").Value = -500514

$ws.Range("H62").Value = 8302.5
$ws.Range("I62").Value = 7095.75
$ws.Range("J62").Value = 9750.6
$ws.Range("K62").Value = 7095.75
$ws.Range("L62").Value = 9750.6
$ws.Range("M62").Value = -6471.75
$ws.Range("N62").Value = -10998.6

$ws.Range("H65").Value = 8302.5
$ws.Range("I65").Value = 7095.75
$ws.Range("J65").Value = 9750.6
$ws.Range("K65").Value = 35478.75
$ws.Range("L65").Value = 48753
$ws.Range("M65").Value = -32358.75
$ws.Range("N65").Value = -54993

$ws.Range("H106").Value = 14002.6
$ws.Range("I106").Value = 2754.3333
$ws.Range("J106").Value = 18823.285
$ws.Range("K106").Value = 2754.3333
$ws.Range("L106").Value = 18823.285
$ws.Range("M106").Value = -2123.3333
$ws.Range("N106").Value = -20085.285

$ws.Range("H112").Value = 1337.9412
$ws.Range("J112").Value = 1496.6
$ws.Range("L112").Value = 4489.799999999999
$ws.Range("N112").Value = -6705.799999999999

$ws.Range("H127").Value = 864.5714
$ws.Range("I127").Value = 546.46155
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 1639.38465
$ws.Range("L127").Value = 15000
$ws.Range("M127").Value = 3320.61535
$ws.Range("N127").Value = -24920

$ws.Range("H138").Value = 2935.6897
$ws.Range("I138").Value = 2143.4546
$ws.Range("J138").Value = 3419.8333
$ws.Range("K138").Value = 6430.3638
$ws.Range("L138").Value = 10259.4999
$ws.Range("M138").Value = -1290.3638
$ws.Range("N138").Value = -20539.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1389
$ws.Range("I4").Value = 185.33333
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 185.33333
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -69.33332999999999
$ws.Range("N4").Value = -5232

$ws.Range("H32").Value = 3635.7
$ws.Range("I32").Value = 3033.4443
$ws.Range("K32").Value = 3033.4443
$ws.Range("M32").Value = -2746.4443

$ws.Range("H61").Value = 8180.6
$ws.Range("I61").Value = 6603.609
$ws.Range("J61").Value = 11203.167
$ws.Range("K61").Value = 6603.609
$ws.Range("L61").Value = 11203.167
$ws.Range("M61").Value = -6391.609
$ws.Range("N61").Value = -11627.167

$ws.Range("H97").Value = 917.5625
$ws.Range("I97").Value = 946.7692
$ws.Range("J97").Value = 791
$ws.Range("K97").Value = 946.7692
$ws.Range("L97").Value = 791
$ws.Range("M97").Value = -450.7692
$ws.Range("N97").Value = -1783

$ws.Range("H136").Value = 8180.6
$ws.Range("I136").Value = 6603.609
$ws.Range("J136").Value = 11203.167
$ws.Range("K136").Value = 19810.827
$ws.Range("L136").Value = 33609.501
$ws.Range("M136").Value = -17260.827
$ws.Range("N136").Value = -38709.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 277.2
$ws.Range("I22").Value = 277.2
$ws.Range("K22").Value = 277.2
$ws.Range("M22").Value = -104.2

$ws.Range("H86").Value = 6322.1
$ws.Range("I86").Value = 3151.875
$ws.Range("J86").Value = 19003
$ws.Range("K86").Value = 3151.875
$ws.Range("L86").Value = 19003
$ws.Range("M86").Value = -2028.875
$ws.Range("N86").Value = -21249

$ws.Range("H89").Value = 6322.1
$ws.Range("I89").Value = 3151.875
$ws.Range("J89").Value = 19003
$ws.Range("K89").Value = 15759.375
$ws.Range("L89").Value = 95015
$ws.Range("M89").Value = -10143.375
$ws.Range("N89").Value = -106247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 991.2857
$ws.Range("I16").Value = 917.8
$ws.Range("J16").Value = 1175
$ws.Range("K16").Value = 917.8
$ws.Range("L16").Value = 1175
$ws.Range("M16").Value = -630.8
$ws.Range("N16").Value = -1749

$ws.Range("H74").Value = 220000

$ws.Range("H77").Value = 220000

$ws.Range("H113").Value = 991.2857
$ws.Range("I113").Value = 917.8
$ws.Range("J113").Value = 1175
$ws.Range("K113").Value = 917.8
$ws.Range("L113").Value = 1175
$ws.Range("M113").Value = 1252.2
$ws.Range("N113").Value = -5515

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 35478
$ws.Range("I2").Value = 168.83333
$ws.Range("J2").Value = 120220
$ws.Range("K2").Value = 1012.99998
$ws.Range("L2").Value = 721320
$ws.Range("M2").Value = -899.9999799999999
$ws.Range("N2").Value = -721546

$ws.Range("H113").Value = 1521.8572
$ws.Range("J113").Value = 2499.5
$ws.Range("L113").Value = 7498.5
$ws.Range("N113").Value = -11838.5

$ws.Range("H119").Value = 33712.332
$ws.Range("I119").Value = 55359.5
$ws.Range("J119").Value = 22888.75
$ws.Range("K119").Value = 166078.5
$ws.Range("L119").Value = 68666.25
$ws.Range("M119").Value = -161240.5
$ws.Range("N119").Value = -78342.25

$ws.Range("H138").Value = 10794
$ws.Range("I138").Value = 9264.286
$ws.Range("K138").Value = 27792.858
$ws.Range("M138").Value = -22652.858

$ws.Range("H139").Value = 4606.636
$ws.Range("I139").Value = 1171.0834
$ws.Range("J139").Value = 8729.299999999999
$ws.Range("K139").Value = 3513.2502
$ws.Range("L139").Value = 26187.9
$ws.Range("M139").Value = 1626.7498
$ws.Range("N139").Value = -36467.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 942
$ws.Range("I2").Value = 115.42857
$ws.Range("K2").Value = 115.42857
$ws.Range("M2").Value = -2.428569999999993

$ws.Range("H35").Value = 18507.5
$ws.Range("I35").Value = 22015
$ws.Range("K35").Value = 22015
$ws.Range("M35").Value = -21717

$ws.Range("H51").Value = 89999.664
$ws.Range("J51").Value = 89999.664
$ws.Range("L51").Value = 89999.664
$ws.Range("N51").Value = -91017.664

$ws.Range("H99").Value = 29685.9
$ws.Range("I99").Value = 16306.667
$ws.Range("K99").Value = 16306.667
$ws.Range("M99").Value = -14060.667

$ws.Range("H126").Value = 9334.666999999999
$ws.Range("I126").Value = 3997.5
$ws.Range("J126").Value = 12003.25
$ws.Range("K126").Value = 11992.5
$ws.Range("L126").Value = 36009.75
$ws.Range("M126").Value = -9522.5
$ws.Range("N126").Value = -40949.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1041.6842
$ws.Range("I16").Value = 1141.2858
$ws.Range("J16").Value = 762.8
$ws.Range("K16").Value = 1141.2858
$ws.Range("L16").Value = 762.8
$ws.Range("M16").Value = -971.2858000000001
$ws.Range("N16").Value = -1102.8

$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5000
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0
$ws.Range("M29").Value = -4705

$ws.Range("H93").Value = 1064.2858
$ws.Range("I93").Value = 1064.2858
$ws.Range("K93").Value = 1064.2858
$ws.Range("M93").Value = 183.7141999999999

$ws.Range("H100").Value = 4694.4116
$ws.Range("I100").Value = 3440.6
$ws.Range("K100").Value = 3440.6
$ws.Range("M100").Value = -2899.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240

$ws.Range("H70").Value = 50000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 50000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H107").Value = 596.8889
$ws.Range("I107").Value = 484.57144
$ws.Range("J107").Value = 990
$ws.Range("K107").Value = 1453.71432
$ws.Range("L107").Value = 2970
$ws.Range("M107").Value = 466.28568
$ws.Range("N107").Value = -6810

$ws.Range("H126").Value = 5367.7
$ws.Range("I126").Value = 5020.8237
$ws.Range("K126").Value = 15062.4711
$ws.Range("M126").Value = -12592.4711
